$wb = $excel.ActiveWorkbook

# --- Sheet "stickers": append two new rows (A3, A4) with "dcdc" ---
$stickers = $wb.Worksheets.Item("stickers")
$stickers.Range("A3").Value = "dcdc"
$stickers.Range("A4").Value = "dcdc"

# --- Sheet "users": append a new row (row 3) with user data ---
$users = $wb.Worksheets.Item("users")
$users.Range("A3").Value = 781523035
$users.Range("B3").Value = "Liran"
$users.Range("C3").Value = "м"
# "11" must be stored as text (like the existing "10н" value), not a number
$users.Range("D3").Value = "'11"
